$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 17, pushing existing rows 17-28 down to 18-29
$ws.Rows.Item(17).Insert()

# Populate the new row 17 with the new data
$ws.Cells.Item(17, 1).Value2 = 1
$ws.Cells.Item(17, 2).Value2 = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(17, 3).Value2 = "Arica y Parinacota"
$ws.Cells.Item(17, 4).Value2 = 44586
$ws.Cells.Item(17, 5).Value2 = 15
$ws.Cells.Item(17, 6).Value2 = 100112028
$ws.Cells.Item(17, 7).Value2 = "Sandia"
$ws.Cells.Item(17, 8).Value2 = "Sin especificar"
$ws.Cells.Item(17, 9).Value2 = "Tercera"
$ws.Cells.Item(17, 10).Value2 = 500
$ws.Cells.Item(17, 11).Value2 = 330
$ws.Cells.Item(17, 12).Value2 = 350
$ws.Cells.Item(17, 13).Value2 = 340
$ws.Cells.Item(17, 14).Value2 = '$/kilo (volumen en unidades)'
$ws.Cells.Item(17, 15).Value2 = "Región de Arica y Parinacota"
$ws.Cells.Item(17, 16).Value2 = 340
$ws.Cells.Item(17, 17).Value2 = 1
$ws.Cells.Item(17, 18).Value2 = "Hortaliza"

Write-Output "done"
